$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 1.042494666666667
$ws.Range("H2").Value = 3.127484
$ws.Range("I2").Value = 0.0007670466909205676
$ws.Range("J2").Value = 0.0007670466909205677
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 63.63554904726622
$ws.Range("R2").Value = 572.7199414253961
$ws.Range("S2").Value = 0.0001567546942454426
$ws.Range("T2").Value = 0.0001567546942454427

# Row 3
$ws.Range("G3").Value = 1.042494666666667
$ws.Range("H3").Value = 3.127484
$ws.Range("I3").Value = 0.0007670466909205676
$ws.Range("J3").Value = 0.0007670466909205677
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("Q3").Value = 110.8322637945146
$ws.Range("R3").Value = 997.490374150632
$ws.Range("S3").Value = 0.0002730152860115182
$ws.Range("T3").Value = 0.0002730152860115183

# Row 4
$ws.Range("G4").Value = 1.042494666666667
$ws.Range("H4").Value = 3.127484
$ws.Range("I4").Value = 0.0007670466909205676
$ws.Range("J4").Value = 0.0007670466909205677
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 136.9195912584835
$ws.Range("R4").Value = 1232.276321326352
$ws.Range("S4").Value = 0.0003372767106636067
$ws.Range("T4").Value = 0.0003372767106636068

# Row 5
$ws.Range("I5").Value = 0.9658609009611662
$ws.Range("J5").Value = 0.9658609009611662
$ws.Range("M5").Value = 61.04160633333334
$ws.Range("N5").Value = 183.124819
$ws.Range("O5").Value = 0.2043613460574534
$ws.Range("P5").Value = 0.2043613460574534
$ws.Range("Q5").Value = 80129.52726800293
$ws.Range("R5").Value = 721165.7454120263
$ws.Range("S5").Value = 0.1973846338246886
$ws.Range("T5").Value = 0.1973846338246886

# Row 6
$ws.Range("I6").Value = 0.9658609009611662
$ws.Range("J6").Value = 0.9658609009611662
$ws.Range("O6").Value = 0.3559304658284363
$ws.Range("P6").Value = 0.3559304658284363
$ws.Range("S6").Value = 0.343779320404581
$ws.Range("T6").Value = 0.3437793204045811

# Row 7
$ws.Range("I7").Value = 0.9658609009611662
$ws.Range("J7").Value = 0.9658609009611662
$ws.Range("M7").Value = 131.3384093333333
$ws.Range("N7").Value = 394.015228
$ws.Range("O7").Value = 0.4397081881141102
$ws.Range("P7").Value = 0.4397081881141103
$ws.Range("Q7").Value = 172408.3831369378
$ws.Range("R7").Value = 1551675.44823244
$ws.Range("S7").Value = 0.4246969467318965
$ws.Range("T7").Value = 0.4246969467318965

# Row 8
$ws.Range("G8").Value = 45.356022
$ws.Range("H8").Value = 136.068066
$ws.Range("I8").Value = 0.03337205234791334
$ws.Range("J8").Value = 0.03337205234791334
$ws.Range("M8").Value = 61.04160633333334
$ws.Range("N8").Value = 183.124819
$ws.Range("O8").Value = 0.2043613460574534
$ws.Range("P8").Value = 0.2043613460574534
$ws.Range("Q8").Value = 2768.604439770006
$ws.Range("R8").Value = 24917.43995793006
$ws.Range("S8").Value = 0.006819957538519369
$ws.Range("T8").Value = 0.006819957538519369

# Row 9
$ws.Range("G9").Value = 45.356022
$ws.Range("H9").Value = 136.068066
$ws.Range("I9").Value = 0.03337205234791334
$ws.Range("J9").Value = 0.03337205234791334
$ws.Range("O9").Value = 0.3559304658284363
$ws.Range("P9").Value = 0.3559304658284363
$ws.Range("Q9").Value = 4822.001258814252
$ws.Range("R9").Value = 43398.01132932826
$ws.Range("S9").Value = 0.01187813013784375
$ws.Range("T9").Value = 0.01187813013784376

# Row 10
$ws.Range("G10").Value = 45.356022
$ws.Range("H10").Value = 136.068066
$ws.Range("I10").Value = 0.03337205234791334
$ws.Range("J10").Value = 0.03337205234791334
$ws.Range("M10").Value = 131.3384093333333
$ws.Range("N10").Value = 394.015228
$ws.Range("O10").Value = 0.4397081881141102
$ws.Range("P10").Value = 0.4397081881141103
$ws.Range("Q10").Value = 5956.98778316767
$ws.Range("R10").Value = 53612.89004850904
$ws.Range("S10").Value = 0.01467396467155021
$ws.Range("T10").Value = 0.01467396467155021